$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D (Supplies stock) for rows 2-6: 2 -> 4
$ws.Range("D2:D6").Value = 4

# Row 5: F5 changes from 4 to 0
$ws.Range("F5").Value = 0

# Row 6: reset to the new plain values (was C6=6, D6=2->4, E6=4->0, F6=0, G6=6->0)
$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = 0

# New row 7
$ws.Range("A7").Value = "-"
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 0

# New row 8
$ws.Range("A8").Value = "-"
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 8
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 8

# Apply the same style as B2:B6 (s="3") to B7 and B8
$ws.Range("B6").Copy()
$ws.Range("B7:B8").PasteSpecial(-4122)  # xlPasteFormats
